$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.472.99"
$ws.Range("E2").Value = "  -2.67%  "

$ws.Range("D3").Value = "3.098.40"
$ws.Range("E3").Value = "  -1.38%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.64"
$ws.Range("E5").Value = "  -2.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.54"
$ws.Range("E6").Value = "  -7.21%  "

$ws.Range("D8").Value = "3.086.72"
$ws.Range("E8").Value = "  -1.50%  "

$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.60"
$ws.Range("E10").Value = "  -4.55%  "

$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.07"
$ws.Range("E13").Value = "  -3.16%  "

$ws.Range("E14").Value = "  -1.21%  "

$ws.Range("D15").Value = "3.602.69"
$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("D16").Value = "63.554.54"
$ws.Range("E16").Value = "  -2.51%  "

$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").Value = "3.096.93"
$ws.Range("E18").Value = "  -1.49%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "489.37"
$ws.Range("E20").Value = "  -6.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.48"
$ws.Range("E21").Value = "  -2.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.704"
$ws.Range("E22").Value = "  +0.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.21"
$ws.Range("E23").Value = "  -2.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.92"
$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("E25").Value = "  -3.26%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  -1.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.33"
$ws.Range("E28").Value = "  -2.76%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.95"
$ws.Range("E30").Value = "  -8.40%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.55"
$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.52"
$ws.Range("E33").Value = "  -6.21%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "59.80"
$ws.Range("E34").Value = "  +12.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "532.00"
$ws.Range("E35").Value = "  -5.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.98"
$ws.Range("E36").Value = "  -1.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.20"
$ws.Range("E37").Value = "  -4.56%  "

$ws.Range("E38").Value = "  -8.35%  "

$ws.Range("E39").Value = "  -2.95%  "

$ws.Range("D40").Value = "3.065.18"
$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  -1.95%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.68"
$ws.Range("E42").Value = "  -6.89%  "

$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.13"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.259"
$ws.Range("E44").Value = "  +0.80%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  -5.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.80"
$ws.Range("E47").Value = "  +3.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.46"
$ws.Range("E48").Value = "  -2.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.107"
$ws.Range("E49").Value = "  -0.90%  "

$ws.Range("D50").Value = "0.0₃0507"
$ws.Range("E50").Value = "  -3.52%  "

$ws.Range("E51").Value = "  +62.37%  "
